$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A12").Value = "IM2325000686"
